$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row groups: rows 2,4,6 share one set of values; rows 3,5 share another set.
$groupA = @(2, 4, 6)
$groupB = @(3, 5)

foreach ($r in $groupA) {
    $ws.Cells.Item($r, 5).Value  = 1000               # E: Max num epochs
    $ws.Cells.Item($r, 11).Value = 1                  # K: Shuffle seed
    $ws.Cells.Item($r, 12).Value = 0.9902999997138977 # L: Accuracy Validate Best
    $ws.Cells.Item($r, 13).Value = 0.9814000129699707 # M: Accuracies Product
    $ws.Cells.Item($r, 14).Value = 432.465            # N: Train time
    $ws.Cells.Item($r, 15).Value = 0.0023             # O: Accuracy Validate per Time
    $ws.Cells.Item($r, 16).Value = 0.0023             # P: Accuracies Product per Time
    $ws.Cells.Item($r, 17).Value = 54                 # Q: Num epochs
    $ws.Cells.Item($r, 18).Value = 8.008599999999999  # R: Average epoch time
    $ws.Cells.Item($r, 19).Value = 0.9901999831199646 # S: Accuracy Validate Last
    $ws.Cells.Item($r, 20).Value = 0.991100013256073  # T: Accuracy Train Last
    $ws.Cells.Item($r, 21).Value = 0.991100013256073  # U: Accuracy Train Best
}

foreach ($r in $groupB) {
    $ws.Cells.Item($r, 5).Value  = 1000               # E: Max num epochs
    $ws.Cells.Item($r, 11).Value = 3                  # K: Shuffle seed
    $ws.Cells.Item($r, 12).Value = 0.9980000257492065 # L: Accuracy Validate Best
    $ws.Cells.Item($r, 13).Value = 0.9959999918937683 # M: Accuracies Product
    $ws.Cells.Item($r, 14).Value = 631.567            # N: Train time
    $ws.Cells.Item($r, 15).Value = 0.0016             # O: Accuracy Validate per Time
    $ws.Cells.Item($r, 16).Value = 0.0016             # P: Accuracies Product per Time
    $ws.Cells.Item($r, 17).Value = 80                 # Q: Num epochs
    $ws.Cells.Item($r, 18).Value = 7.8946             # R: Average epoch time
    $ws.Cells.Item($r, 19).Value = 0.9977999925613403 # S: Accuracy Validate Last
    $ws.Cells.Item($r, 20).Value = 0.9980000257492065 # T: Accuracy Train Last
    $ws.Cells.Item($r, 21).Value = 0.9980000257492065 # U: Accuracy Train Best
}
